# Append the new profit row (2025-11-06) that the diff adds as row 81.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as literal text ("MM/DD/YYYY"), matching the
# existing rows in the sheet. Force text formatting first so Excel doesn't
# auto-convert the string into a date serial number, then clear the
# formatting override again so the cell ends up unstyled like its siblings.
$ws.Cells.Item(81, 1).NumberFormat = "@"
$ws.Cells.Item(81, 1).Value = "11/06/2025"
$ws.Cells.Item(81, 1).ClearFormats()

# Column B holds the numeric profit value.
$ws.Cells.Item(81, 2).Value = 8915.879999999999
